$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5's mobile number should become a true numeric value instead of text
$ws.Range("D5").Value = 9386776366

# New row 6: Abhinab Kuamr / 89871BZ / Trouble opening the portal
$ws.Range("A6").Value = "2025-06-04T23:11:42.841563"
$ws.Range("B6").Value = "Abhinab Kuamr"
$ws.Range("C6").Value = "89871BZ"
$ws.Range("D6").Value = "'"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "Trouble opening the portal"

# New row 7: Manish Guota / 8724B / unable to open loan account...
$ws.Range("A7").Value = "2025-06-04T23:12:15.948330"
$ws.Range("B7").Value = "Manish Guota"
$ws.Range("C7").Value = "8724B"
$ws.Range("D7").Value = "'"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "unable to open loan account. Name: Manish Guota, id : 8724B"
